$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '43.213.61'
$ws.Cells.Item(2, 5).Value = '  +0.95%  '
$ws.Cells.Item(3, 4).Value = '2.377.31'
$ws.Cells.Item(3, 5).Value = '  +3.11%  '
$ws.Cells.Item(4, 5).Value = '  -0.06%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '303.42'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '97.41'
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  +1.37%  '
$ws.Cells.Item(7, 5).Value = '  -0.03%  '
$ws.Cells.Item(8, 5).Value = '  -0.12%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.502'
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  +1.62%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '34.20'
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -1.22%  '
$ws.Cells.Item(11, 5).Value = '  +0.26%  '
$ws.Cells.Item(12, 5).Value = '  +3.08%  '
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '18.54'
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  -3.38%  '
$ws.Cells.Item(14, 5).Value = '  -0.04%  '
$ws.Cells.Item(15, 4).Value = '2.744.33'
$ws.Cells.Item(15, 5).Value = '  +2.91%  '
$ws.Cells.Item(16, 4).Value = '2.367.38'
$ws.Cells.Item(16, 5).Value = '  +2.46%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '0.807'
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  +3.02%  '
$ws.Cells.Item(18, 4).Value = '43.194.98'
$ws.Cells.Item(18, 5).Value = '  +1.06%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '12.32'
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  +1.09%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '6.33'
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  +5.29%  '
$ws.Cells.Item(21, 5).Value = '  -0.17%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '68.30'
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +0.78%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '236.16'
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +0.28%  '
$ws.Cells.Item(24, 5).Value = '  -2.00%  '
$ws.Cells.Item(25, 5).Value = '  +0.81%  '
$ws.Cells.Item(26, 5).Value = '  +0.01%  '
$ws.Cells.Item(27, 5).Value = '  +1.75%  '
$ws.Cells.Item(28, 5).Value = '  +0.09%  '
$ws.Cells.Item(29, 5).Value = '  +0.98%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '31.63'
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -1.35%  '
$ws.Cells.Item(31, 2).Value = 'Filecoin'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '5.13'
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  +2.84%  '
$ws.Cells.Item(32, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '1.00'
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -0.06%  '
$ws.Cells.Item(33, 5).Value = '  +3.88%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '17.24'
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  -2.05%  '
$ws.Cells.Item(35, 5).Value = '  +5.85%  '
$ws.Cells.Item(37, 5).Value = '  -1.04%  '
$ws.Cells.Item(38, 2).Value = 'EnergySwap'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '22.91'
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  +13.18%  '
$ws.Cells.Item(39, 2).Value = 'Kaspa'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.101'
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  +1.30%  '
$ws.Cells.Item(40, 5).Value = '  +3.50%  '
$ws.Cells.Item(41, 5).Value = '  +0.17%  '
$ws.Cells.Item(42, 4).Value = '1.948.41'
$ws.Cells.Item(42, 5).Value = '  -0.95%  '
$ws.Cells.Item(43, 2).Value = 'Monero'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '101.85'
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  -38.50%  '
$ws.Cells.Item(44, 2).Value = 'VeChain'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.0281'
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  +0.60%  '
$ws.Cells.Item(45, 5).Value = '  +3.06%  '
$ws.Cells.Item(46, 5).Value = '  -10.10%  '
$ws.Cells.Item(47, 5).Value = '  -0.40%  '
$ws.Cells.Item(48, 4).Value = '2.599.53'
$ws.Cells.Item(48, 5).Value = '  +2.73%  '
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '53.12'
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -0.31%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '72.31'
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  +1.20%  '
